$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-11 18:29:16"

for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
